# Add data for 2023-10-16: update the "2023" column (J) totals for
# Citywide Totals, By Neighborhood, and each affected per-neighborhood
# sheet with the refreshed year-to-date crime counts.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6074
$ws.Range("J3").Value = 6486
$ws.Range("J4").Value = 1401
$ws.Range("J5").Value = 496
$ws.Range("J6").Value = 8418
$ws.Range("J7").Value = 22875

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 181
$ws.Range("J7").Value = 673
$ws.Range("J8").Value = 1435
$ws.Range("J10").Value = 163
$ws.Range("J11").Value = 367
$ws.Range("J15").Value = 254
$ws.Range("J18").Value = 193
$ws.Range("J19").Value = 677
$ws.Range("J20").Value = 482
$ws.Range("J21").Value = 65
$ws.Range("J23").Value = 213
$ws.Range("J26").Value = 48
$ws.Range("J29").Value = 1258
$ws.Range("J32").Value = 37
$ws.Range("J36").Value = 313
$ws.Range("J37").Value = 698
$ws.Range("J42").Value = 964
$ws.Range("J43").Value = 193
$ws.Range("J44").Value = 176
$ws.Range("J48").Value = 270
$ws.Range("J51").Value = 289
$ws.Range("J52").Value = 573
$ws.Range("J53").Value = 321
$ws.Range("J54").Value = 441
$ws.Range("J55").Value = 327
$ws.Range("J56").Value = 30
$ws.Range("J63").Value = 79
$ws.Range("J64").Value = 150
$ws.Range("J65").Value = 564
$ws.Range("J66").Value = 69
$ws.Range("J67").Value = 868
$ws.Range("J72").Value = 91
$ws.Range("J73").Value = 221
$ws.Range("J76").Value = 351
$ws.Range("J77").Value = 173
$ws.Range("J79").Value = 652
$ws.Range("J83").Value = 456
$ws.Range("J84").Value = 191
$ws.Range("J85").Value = 938
$ws.Range("J89").Value = 301
$ws.Range("J90").Value = 247
$ws.Range("J91").Value = 264
$ws.Range("J93").Value = 100
$ws.Range("J94").Value = 235
$ws.Range("J95").Value = 331
$ws.Range("J97").Value = 205
$ws.Range("J99").Value = 355
$ws.Range("J101").Value = 22875

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 206
$ws.Range("J7").Value = 673

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 159
$ws.Range("J7").Value = 367

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 86
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 250
$ws.Range("J6").Value = 272
$ws.Range("J7").Value = 938

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J4").Value = 20
$ws.Range("J6").Value = 235
$ws.Range("J7").Value = 573

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 42
$ws.Range("J6").Value = 211
$ws.Range("J7").Value = 321

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 392
$ws.Range("J6").Value = 494
$ws.Range("J7").Value = 1435

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 138
$ws.Range("J3").Value = 170
$ws.Range("J7").Value = 456

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 114
$ws.Range("J7").Value = 331

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 202
$ws.Range("J7").Value = 698

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 165
$ws.Range("J3").Value = 161
$ws.Range("J7").Value = 564

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 329
$ws.Range("J6").Value = 235
$ws.Range("J7").Value = 868

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 212
$ws.Range("J7").Value = 441

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 386
$ws.Range("J6").Value = 322
$ws.Range("J7").Value = 1258

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 270

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 163
$ws.Range("J6").Value = 260
$ws.Range("J7").Value = 677

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 58
$ws.Range("J6").Value = 195
$ws.Range("J7").Value = 351

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 207
$ws.Range("J3").Value = 196
$ws.Range("J6").Value = 501
$ws.Range("J7").Value = 964

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 32
$ws.Range("J6").Value = 88
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 70
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 171
$ws.Range("J7").Value = 327

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 72
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J2").Value = 9
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J6").Value = 191
$ws.Range("J7").Value = 652

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 134
$ws.Range("J3").Value = 167
$ws.Range("J7").Value = 482

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 101
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 313

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J3").Value = 10
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 143
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 89
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 76
$ws.Range("J5").Value = 8
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 66
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 30
